$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 25.8200000000006
$ws.Range("H2").Value = [double]"3.325564629708122e-05"
$ws.Range("I2").Value = [double]"3.325564629708122e-05"
$ws.Range("L2").Value = 36.88085379122087
$ws.Range("M2").Value = "[20.36872753723, 53.39298004521174]"
$ws.Range("N2").Value = [double]"4.781287169541848e-05"
$ws.Range("O2").Value = [double]"4.781287169541848e-05"
$ws.Range("P2").Value = 1.603816069400194
$ws.Range("Q2").Value = "[1.0126054398958084, 2.19502669890458]"
$ws.Range("R2").Value = [double]"1.936909117761232e-06"
$ws.Range("S2").Value = [double]"1.936909117761232e-06"
$ws.Range("T2").Value = 61.22241615376245
$ws.Range("U2").Value = "[50.573268746141224, 71.87156356138367]"
$ws.Range("V2").Value = [double]"4.218847493575595e-15"
$ws.Range("W2").Value = [double]"4.218847493575595e-15"
$ws.Range("X2").Value = 19.22930930930976
$ws.Range("Y2").Value = 16.79979979980019
$ws.Range("Z2").Value = 21.65881881881932
$ws.Range("F3").Value = 25.8200000000006
$ws.Range("H3").Value = [double]"1.51346122695406e-05"
$ws.Range("I3").Value = [double]"1.51346122695406e-05"
$ws.Range("L3").Value = 45.20441970854078
$ws.Range("M3").Value = "[22.403710615243455, 68.00512880183811]"
$ws.Range("N3").Value = 0.0002383024260534938
$ws.Range("O3").Value = 0.0002383024260534938
$ws.Range("P3").Value = 0.9622896416401163
$ws.Range("Q3").Value = "[0.4465527095192696, 1.478026573760963]"
$ws.Range("R3").Value = 0.0004908064229998121
$ws.Range("S3").Value = 0.0004908064229998121
$ws.Range("T3").Value = 60.13486599302988
$ws.Range("U3").Value = "[48.19213535403905, 72.0775966320207]"
$ws.Range("V3").Value = [double]"3.341771304121721e-13"
$ws.Range("W3").Value = [double]"3.341771304121721e-13"
$ws.Range("X3").Value = 21.86558558558609
$ws.Range("Y3").Value = 19.74622622622669
$ws.Range("Z3").Value = 23.9849449449455
$ws.Range("F4").Value = 25.8200000000006
$ws.Range("H4").Value = [double]"4.081925162990885e-06"
$ws.Range("I4").Value = [double]"4.081925162990885e-06"
$ws.Range("L4").Value = 46.00818609867494
$ws.Range("M4").Value = "[24.884215565250685, 67.1321566320992]"
$ws.Range("N4").Value = [double]"6.860163279220899e-05"
$ws.Range("O4").Value = [double]"6.860163279220899e-05"
$ws.Range("P4").Value = 0.5849211547224229
$ws.Range("Q4").Value = "[0.11950002085726741, 1.0503422885875784]"
$ws.Range("R4").Value = 0.01492979866304922
$ws.Range("S4").Value = 0.01492979866304922
$ws.Range("T4").Value = 54.40979872784671
$ws.Range("U4").Value = "[43.370712456137795, 65.44888499955563]"
$ws.Range("V4").Value = [double]"6.534772722943671e-13"
$ws.Range("W4").Value = [double]"6.534772722943671e-13"
$ws.Range("X4").Value = 23.41633633633688
$ws.Range("Y4").Value = 21.50374374374424
$ws.Range("Z4").Value = 25.32892892892952
$ws.Range("F5").Value = 25.8200000000006
$ws.Range("H5").Value = 0.02776208035334948
$ws.Range("I5").Value = 0.02776208035334948
$ws.Range("L5").Value = 26.68082237112299
$ws.Range("M5").Value = "[4.760513507497663, 48.60113123474832]"
$ws.Range("N5").Value = 0.01817273898572802
$ws.Range("O5").Value = 0.01817273898572802
$ws.Range("P5").Value = -0.06289474781961513
$ws.Range("Q5").Value = "[-1.2264475824825016, 1.1006580868432714]"
$ws.Range("R5").Value = 0.9137892698473666
$ws.Range("S5").Value = 0.9137892698473666
$ws.Range("T5").Value = 49.47285913775396
$ws.Range("U5").Value = "[36.22916920186075, 62.716549073647165]"
$ws.Range("V5").Value = [double]"1.714974384725565e-09"
$ws.Range("W5").Value = [double]"1.714974384725565e-09"
$ws.Range("X5").Value = 0.2584584584584633
$ws.Range("Y5").Value = -4.523023023023129
$ws.Range("Z5").Value = 5.039939939940056
$ws.Range("F6").Value = 25.8200000000006
$ws.Range("H6").Value = 0.0002097029021311814
$ws.Range("I6").Value = 0.0002097029021311814
$ws.Range("L6").Value = 34.64304249637465
$ws.Range("M6").Value = "[15.159202190730113, 54.126882802019196]"
$ws.Range("N6").Value = 0.0008349725165663902
$ws.Range("O6").Value = 0.0008349725165663902
$ws.Range("P6").Value = -0.3396316382259235
$ws.Range("Q6").Value = "[-0.9811580659860013, 0.3018947895341544]"
$ws.Range("R6").Value = 0.291980167510081
$ws.Range("S6").Value = 0.291980167510081
$ws.Range("T6").Value = 48.11191783498744
$ws.Range("U6").Value = "[37.34967002417239, 58.874165645802485]"
$ws.Range("V6").Value = [double]"1.251065917529104e-11"
$ws.Range("W6").Value = [double]"1.251065917529104e-11"
$ws.Range("X6").Value = 1.395675675675704
$ws.Range("Y6").Value = -1.240600600600632
$ws.Range("Z6").Value = 4.031951951952042
$ws.Range("F7").Value = 25.8200000000006
$ws.Range("H7").Value = [double]"4.673062797255767e-05"
$ws.Range("I7").Value = [double]"4.673062797255767e-05"
$ws.Range("L7").Value = 53.52385171869542
$ws.Range("M7").Value = "[24.31140409872492, 82.73629933866592]"
$ws.Range("N7").Value = 0.0006023353738586756
$ws.Range("O7").Value = 0.0006023353738586756
$ws.Range("P7").Value = -0.4780000834290776
$ws.Range("Q7").Value = "[-1.0440528138056173, 0.08805264694746207]"
$ws.Range("R7").Value = 0.09588264994104545
$ws.Range("S7").Value = 0.09588264994104545
$ws.Range("T7").Value = 69.80716225190571
$ws.Range("U7").Value = "[54.5703373718108, 85.04398713200061]"
$ws.Range("V7").Value = [double]"6.063594071292755e-12"
$ws.Range("W7").Value = [double]"6.063594071292755e-12"
$ws.Range("X7").Value = 1.964284284284329
$ws.Range("Y7").Value = -0.3618418418418528
$ws.Range("Z7").Value = 4.290410410410511
$ws.Range("F8").Value = 25.8200000000006
$ws.Range("H8").Value = [double]"7.9558035492866e-08"
$ws.Range("I8").Value = [double]"7.9558035492866e-08"
$ws.Range("L8").Value = 62.80320777866051
$ws.Range("M8").Value = "[38.176355493121406, 87.43006006419962]"
$ws.Range("N8").Value = [double]"5.834441962671022e-06"
$ws.Range("O8").Value = [double]"5.834441962671022e-06"
$ws.Range("P8").Value = -1.182421259008771
$ws.Range("Q8").Value = "[-1.6101055441821552, -0.754736973835386]"
$ws.Range("R8").Value = [double]"1.358742622192821e-06"
$ws.Range("S8").Value = [double]"1.358742622192821e-06"
$ws.Range("T8").Value = 68.94099619271289
$ws.Range("U8").Value = "[55.513598480450526, 82.36839390497525]"
$ws.Range("V8").Value = [double]"1.798561299892754e-13"
$ws.Range("W8").Value = [double]"1.798561299892754e-13"
$ws.Range("X8").Value = 4.859019019019129
$ws.Range("Y8").Value = 3.101501501501573
$ws.Range("Z8").Value = 6.616536536536685
$ws.Range("F9").Value = 22
$ws.Range("H9").Value = [double]"7.378257696688539e-05"
$ws.Range("I9").Value = [double]"7.378257696688539e-05"
$ws.Range("J9").Value = 0.3359349050182568
$ws.Range("K9").Value = 0.3359349050182568
$ws.Range("L9").Value = 47.16825568737999
$ws.Range("M9").Value = "[21.10564074256405, 73.23087063219593]"
$ws.Range("N9").Value = 0.0006898461948472434
$ws.Range("O9").Value = 0.0006898461948472434
$ws.Range("P9").Value = -2.138421425866927
$ws.Range("Q9").Value = "[-2.6793162571156195, -1.597526594618234]"
$ws.Range("R9").Value = [double]"3.90576682107735e-10"
$ws.Range("S9").Value = [double]"3.90576682107735e-10"
$ws.Range("T9").Value = 62.61421707969022
$ws.Range("U9").Value = "[49.17460349456814, 76.0538306648123]"
$ws.Range("V9").Value = [double]"3.671063453225543e-12"
$ws.Range("W9").Value = [double]"3.671063453225543e-12"
$ws.Range("X9").Value = 7.487487487487488
$ws.Range("Y9").Value = 5.593593593593596
$ws.Range("Z9").Value = 9.381381381381381
$ws.Range("F10").Value = 22
$ws.Range("H10").Value = 0.0006019168706263844
$ws.Range("I10").Value = 0.0006019168706263844
$ws.Range("J10").Value = 0.8108801066811639
$ws.Range("K10").Value = 0.8108801066811639
$ws.Range("L10").Value = 35.59368118719019
$ws.Range("M10").Value = "[14.509301126155655, 56.678061248224715]"
$ws.Range("N10").Value = 0.001421337968732672
$ws.Range("O10").Value = 0.001421337968732672
$ws.Range("P10").Value = -1.912000333716311
$ws.Range("Q10").Value = "[-2.591263610168158, -1.2327370572644631]"
$ws.Range("R10").Value = [double]"9.643522058677689e-07"
$ws.Range("S10").Value = [double]"9.643522058677689e-07"
$ws.Range("T10").Value = 55.81571493586053
$ws.Range("U10").Value = "[44.1313705389788, 67.50005933274227]"
$ws.Range("V10").Value = [double]"1.718403197514817e-12"
$ws.Range("W10").Value = [double]"1.718403197514817e-12"
$ws.Range("X10").Value = 6.694694694694693
$ws.Range("Y10").Value = 4.316316316316314
$ws.Range("Z10").Value = 9.073073073073072
$ws.Range("F11").Value = 22
$ws.Range("H11").Value = 0.005278458173854417
$ws.Range("I11").Value = 0.005278458173854417
$ws.Range("J11").Value = 0.3179276758449368
$ws.Range("K11").Value = 0.3179276758449368
$ws.Range("L11").Value = 35.12705834910162
$ws.Range("M11").Value = "[7.348879767325407, 62.90523693087783]"
$ws.Range("N11").Value = 0.0143566300992497
$ws.Range("O11").Value = 0.0143566300992497
$ws.Range("P11").Value = -1.018894914677771
$ws.Range("Q11").Value = "[-1.7736318885131572, -0.2641579408423853]"
$ws.Range("R11").Value = 0.009267905898218354
$ws.Range("S11").Value = 0.009267905898218354
$ws.Range("T11").Value = 58.72717488280274
$ws.Range("U11").Value = "[44.16606484597425, 73.28828491963124]"
$ws.Range("V11").Value = [double]"2.282791733421163e-10"
$ws.Range("W11").Value = [double]"2.282791733421163e-10"
$ws.Range("X11").Value = 3.567567567567568
$ws.Range("Y11").Value = 0.9249249249249245
$ws.Range("Z11").Value = 6.210210210210212
